# spring 24 week 9 inputs
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New matchup data rows to append starting at row 1508
$data = @(
    @(1,6,3,14),
    @(3,13,4,7),
    @(8,14,6,6),
    @(3,15,2,5),
    @(4,17,2,3),
    @(2,1,4,19),
    @(4,5,3,15),
    @(7,15,6,5),
    @(3,6,2,14),
    @(2,12,4,8),
    @(4,8,3,12),
    @(5,8,4,12),
    @(5,13,6,7),
    @(2,7,5,13),
    @(4,12,3,8),
    @(4,12,5,8)
)

$startRow = 1508
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}

$lastRow = $startRow + $data.Count - 1

# Update the view so the newly added rows are visible/selected, matching the authored edit
$excel.ActiveWindow.ScrollRow = $lastRow - 18
$ws.Range("A" + ($lastRow + 1)).Select()
